# Updates the "Price" (column D) and "Volume(1h)" (column E) values on the
# cryptos worksheet to reflect the latest GitHub Actions refresh.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row number -> new Price (column D) / Volume(1h) (column E) values.
# Only the cells that actually changed are listed; rows/cells not present
# here are left untouched.
$updates = @(
    @{ Row = 2;  D = "36.630.38";    E = "  -0.81%  " },
    @{ Row = 3;  D = "2.062.75";     E = "  +1.22%  " },
    @{ Row = 4;  E = "  +0.02%  " },
    @{ Row = 5;  D = "243.80";       E = "  -0.45%  " },
    @{ Row = 6;  D = "0.669";        E = "  +1.82%  " },
    @{ Row = 7;  E = "  +0.03%  " },
    @{ Row = 8;  D = "54.80";        E = "  -6.53%  " },
    @{ Row = 9;  D = "58.84";        E = "  -0.69%  " },
    @{ Row = 10; E = "  -3.25%  " },
    @{ Row = 11; D = "0.0753" },
    @{ Row = 12; E = "  -3.09%  " },
    @{ Row = 13; D = "0.937";        E = "  +6.99%  " },
    @{ Row = 14; E = "  -3.29%  " },
    @{ Row = 15; D = "2.363.24";     E = "  +1.19%  " },
    @{ Row = 16; E = "  -2.50%  " },
    @{ Row = 17; D = "2.063.43";     E = "  +1.62%  " },
    @{ Row = 18; D = "36.578.98";    E = "  -0.88%  " },
    @{ Row = 19; D = "16.98";        E = "  -6.49%  " },
    @{ Row = 20; D = "72.11";        E = "  -1.90%  " },
    @{ Row = 21; D = "0.0₃0863";     E = "  -2.49%  " },
    @{ Row = 22; D = "238.62";       E = "  +1.46%  " },
    @{ Row = 23; E = "  -1.26%  " },
    @{ Row = 24; E = "  +0.05%  " },
    @{ Row = 25; D = "2.36";         E = "  -3.69%  " },
    @{ Row = 26; D = "2.15";         E = "  +2.11%  " },
    @{ Row = 27; D = "9.35";         E = "  -2.53%  " },
    @{ Row = 28; D = "164.88";       E = "  -1.98%  " },
    @{ Row = 29; D = "20.19";        E = "  +1.45%  " },
    @{ Row = 30; D = "0.122";        E = "  -1.03%  " },
    @{ Row = 31; E = "  +8.62%  " },
    @{ Row = 32; D = "5.10";         E = "  -6.29%  " },
    @{ Row = 33; E = "  -4.32%  " },
    @{ Row = 34; D = "0.0600";       E = "  -1.91%  " },
    @{ Row = 35; E = "  +0.04%  " },
    @{ Row = 36; E = "  -0.47%  " },
    @{ Row = 37; E = "  -0.64%  " },
    @{ Row = 38; D = "0.0824";       E = "  -4.55%  " },
    @{ Row = 39; D = "1.26";         E = "  -3.20%  " },
    @{ Row = 40; D = "4.86";         E = "  -5.71%  " },
    @{ Row = 41; E = "  -1.93%  " },
    @{ Row = 42; D = "2.88";         E = "  -7.61%  " },
    @{ Row = 43; E = "  -2.10%  " },
    @{ Row = 44; D = "94.60";        E = "  -2.47%  " },
    @{ Row = 45; D = "0.0913";       E = "  -4.92%  " },
    @{ Row = 46; D = "1.413.45";     E = "  +9.35%  " },
    @{ Row = 47; D = "16.06";        E = "  -4.43%  " },
    @{ Row = 48; D = "7.51";         E = "  +12.33%  " },
    @{ Row = 49; E = "  +0.57%  " },
    @{ Row = 50; D = "2.28" },
    @{ Row = 51; D = "2.248.15";     E = "  +1.19%  " }
)

foreach ($u in $updates) {
    $row = $u.Row

    if ($u.ContainsKey("D")) {
        $cell = $ws.Range("D$row")
        # Some refreshed prices (e.g. "243.80", "16.98") are valid numeric
        # literals; writing them straight to .Value would make Excel coerce
        # them into numbers. Force text entry via a temporary Text number
        # format, then strip the format again so the cell is left exactly as
        # it started (General, unstyled) but still holds the string value.
        $cell.NumberFormat = "@"
        $cell.Value = $u.D
        $cell.ClearFormats()
    }

    if ($u.ContainsKey("E")) {
        $ws.Range("E$row").Value = $u.E
    }
}
